{"js": "// \"copy 4 links in the song\" \u2014 duplicate the 4-paragraph stanza\n//   \u0418 \u043a\u043e\u043c\u043f\u043e\u0442 \u043d\u0435 \u043b\u044c\u0435\u0442\u044f \u0432 \u0440\u043e\u0442.\u041d\u043e\u0447\u044c\u044e \u0441\u0442\u0430\u043d\u0443 \u043e \u043e\u043a\u043d\u0430\n//   \u0418 \u0441\u0442\u043e\u044e \u0432\u0441\u044e \u043d\u043e\u0447\u044c \u0431\u0435\u0437 \u0441\u043d\u0430,\n//   \u0410 \u0432\u0441\u0451 \u0432\u043e\u043b\u043d\u0443\u044e\u0441\u044c \u043e\u0431 \u0420\u0430\u0441\u0441\u0435\u0435\n//   \u00ab\u041a\u0430\u043a \u0442\u0430\u043c \u0431\u0435\u0434\u043d\u0430\u044f \u043e\u043d\u0430?\u00bb\n// immediately after the paragraph \"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\" (i.e. right\n// before the existing copy of that same stanza), so it ends up appearing\n// twice in a row.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\") by its text.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error('Could not find anchor paragraph \"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\"');\n}\n\n// The OOXML for the 4 paragraphs being duplicated (same markup, including\n// the proofErr spell/grammar-check hints, as the existing copy).\nconst stanzaOoxml =\n  '<w:p><w:r><w:t xml:space=\"preserve\">\u0418 \u043a\u043e\u043c\u043f\u043e\u0442 \u043d\u0435 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u043b\u044c\u0435\u0442\u044f</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u0432 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>\u0440\u043e\u0442.\u041d\u043e\u0447\u044c\u044e</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u0441\u0442\u0430\u043d\u0443 \u043e \u043e\u043a\u043d\u0430</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>\u0418 \u0441\u0442\u043e\u044e \u0432\u0441\u044e \u043d\u043e\u0447\u044c \u0431\u0435\u0437 \u0441\u043d\u0430,</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">\u0410 \u0432\u0441\u0451 \u0432\u043e\u043b\u043d\u0443\u044e\u0441\u044c \u043e\u0431 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u0420\u0430\u0441\u0441\u0435\u0435</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>' +\n  '<w:p><w:r><w:t>\u00ab\u041a\u0430\u043a \u0442\u0430\u043c \u0431\u0435\u0434\u043d\u0430\u044f \u043e\u043d\u0430?\u00bb</w:t></w:r></w:p>';\n\nconst packageOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + stanzaOoxml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nconst insertionRange = anchor.getRange(\"After\");\ninsertionRange.insertOoxml(packageOoxml, \"After\");\nawait context.sync();\n", "ps1": "# \"copy 4 links in the song\" \u2014 duplicate the 4-paragraph stanza\n#   \u0418 \u043a\u043e\u043c\u043f\u043e\u0442 \u043d\u0435 \u043b\u044c\u0435\u0442\u044f \u0432 \u0440\u043e\u0442.\u041d\u043e\u0447\u044c\u044e \u0441\u0442\u0430\u043d\u0443 \u043e \u043e\u043a\u043d\u0430\n#   \u0418 \u0441\u0442\u043e\u044e \u0432\u0441\u044e \u043d\u043e\u0447\u044c \u0431\u0435\u0437 \u0441\u043d\u0430,\n#   \u0410 \u0432\u0441\u0451 \u0432\u043e\u043b\u043d\u0443\u044e\u0441\u044c \u043e\u0431 \u0420\u0430\u0441\u0441\u0435\u0435\n#   \u00ab\u041a\u0430\u043a \u0442\u0430\u043c \u0431\u0435\u0434\u043d\u0430\u044f \u043e\u043d\u0430?\u00bb\n# immediately after the paragraph \"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\" (i.e. right\n# before the existing copy of that same stanza), so it ends up appearing\n# twice in a row.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\") with Find.\n$searchRng = $d.Content\n$found = $searchRng.Find.Execute(\"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\")\nif (-not $found) {\n    throw 'Could not find anchor paragraph \"\u0418 \u0438\u043a\u0440\u0430 \u043d\u0435 \u043b\u0435\u0437\u0435\u0442 \u0432 \u0433\u043e\u0440\u043b\u043e,\"'\n}\n\n# Collapse to the point right after the matched text (before its paragraph\n# mark) and re-materialize as a plain Range so the subsequent InsertXML\n# lands exactly between the anchor paragraph and the stanza that follows it.\n$searchRng.Collapse(0)\n$insertionPoint = $searchRng.Start\n$insertRng = $d.Range($insertionPoint, $insertionPoint)\n\n# OOXML for the 4 paragraphs being duplicated (same markup, including the\n# proofErr spell/grammar-check hints, as the existing copy).\n$stanza = '<w:p><w:r><w:t xml:space=\"preserve\">\u0418 \u043a\u043e\u043c\u043f\u043e\u0442 \u043d\u0435 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u043b\u044c\u0435\u0442\u044f</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u0432 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>\u0440\u043e\u0442.\u041d\u043e\u0447\u044c\u044e</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u0441\u0442\u0430\u043d\u0443 \u043e \u043e\u043a\u043d\u0430</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>\u0418 \u0441\u0442\u043e\u044e \u0432\u0441\u044e \u043d\u043e\u0447\u044c \u0431\u0435\u0437 \u0441\u043d\u0430,</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">\u0410 \u0432\u0441\u0451 \u0432\u043e\u043b\u043d\u0443\u044e\u0441\u044c \u043e\u0431 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>\u0420\u0430\u0441\u0441\u0435\u0435</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>' +\n  '<w:p><w:r><w:t>\u00ab\u041a\u0430\u043a \u0442\u0430\u043c \u0431\u0435\u0434\u043d\u0430\u044f \u043e\u043d\u0430?\u00bb</w:t></w:r></w:p>'\n\n$packageXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $stanza + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$insertRng.InsertXML($packageXml)\n"}
